# fall 23 week 7 updates
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Nine")

# Row 3
$ws.Range("D3").Value = 10.53
$ws.Range("F3").Value = 9.56
$ws.Range("H3").Value = 12.8

# Row 4
$ws.Range("C4").Value = 9.470000000000001
$ws.Range("E4").Value = 10.24
$ws.Range("F4").Value = 10.34
$ws.Range("G4").Value = 10.38

# Row 5
$ws.Range("D5").Value = 9.76
$ws.Range("F5").Value = 10.21
$ws.Range("G5").Value = 9.27

# Row 6
$ws.Range("C6").Value = 10.44
$ws.Range("D6").Value = 9.66
$ws.Range("E6").Value = 9.789999999999999

# Row 7
$ws.Range("D7").Value = 9.619999999999999
$ws.Range("E7").Value = 10.73

# Row 8
$ws.Range("C8").Value = 7.2
